$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
Write-Host $ws.Name
$ws.Range("H15").Value = 4655.04
Write-Host $ws.Range("H15").Value
